$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Price" (column D) values. These strings must be written verbatim
# (several look like numbers, e.g. "7.190" or "0.8680") so force the cells to
# Text format first to stop Excel from re-parsing/truncating them as numbers.
$priceUpdates = @{
    2 = "26.295.06"
    3 = "1.664.92"
    5 = "218.97"
    6 = "0.5354"
    8 = "0.2660"
    9 = "0.06404"
    10 = "20.59"
    11 = "0.07821"
    12 = "4.566"
    13 = "1.669.42"
    14 = "1.892.11"
    15 = "0.5525"
    16 = "0.0₅8205"
    17 = "65.73"
    19 = "4.682"
    20 = "193.75"
    21 = "10.22"
    22 = "6.043"
    24 = "146.14"
    25 = "0.1236"
    26 = "7.190"
    27 = "16.11"
    28 = "1.484"
    29 = "0.05854"
    30 = "1.285"
    31 = "3.611"
    32 = "3.282"
    34 = "0.9646"
    35 = "2.827"
    36 = "2.418"
    38 = "0.01607"
    39 = "0.8680"
    40 = "5.877"
    41 = "1.053.44"
    43 = "104.89"
    44 = "1.802.70"
    45 = "57.87"
    47 = "1.013"
    48 = "0.4384"
    49 = "7.996"
    50 = "0.05163"
}

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$row]
}

# Updated "Volume(1h)" (column E) values (percent strings padded with spaces).
$volumeUpdates = @{
    2 = "  +0.78%  "
    3 = "  +0.59%  "
    4 = "  +0.80%  "
    5 = "  +0.68%  "
    6 = "  +1.84%  "
    7 = "  +0.76%  "
    8 = "  +1.91%  "
    9 = "  +0.91%  "
    10 = "  +0.95%  "
    11 = "  +0.35%  "
    12 = "  +1.39%  "
    13 = "  +0.90%  "
    14 = "  +0.52%  "
    15 = "  +0.69%  "
    16 = "  -0.22%  "
    17 = "  +0.50%  "
    18 = "  +0.87%  "
    19 = "  +2.17%  "
    20 = "  +1.60%  "
    21 = "  +1.58%  "
    22 = "  +0.21%  "
    23 = "  +0.80%  "
    24 = "  +3.04%  "
    25 = "  +0.19%  "
    26 = "  -0.64%  "
    27 = "  +0.25%  "
    28 = "  +3.77%  "
    29 = "  -0.25%  "
    30 = "  +0.96%  "
    31 = "  +2.40%  "
    32 = "  +0.57%  "
    33 = "  +1.62%  "
    34 = "  +1.62%  "
    35 = "  +1.56%  "
    36 = "  +0.35%  "
    37 = "  +1.86%  "
    38 = "  -0.44%  "
    39 = "  +2.48%  "
    40 = "  +1.55%  "
    41 = "  +2.68%  "
    42 = "  +0.75%  "
    43 = "  +2.08%  "
    44 = "  +0.26%  "
    45 = "  +1.44%  "
    46 = "  -5.50%  "
    48 = "  +1.99%  "
    49 = "  +2.57%  "
    50 = "  +0.23%  "
    51 = "  -3.81%  "
}

foreach ($row in $volumeUpdates.Keys) {
    $ws.Range("E$row").Value = $volumeUpdates[$row]
}
